$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").Value = "units"
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 15).Value = "None"
}

$ws.Range("O1:O29").HorizontalAlignment = -4108

$ws.Range("O34").Select()
